$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Index", "Time", "SegStart", "SegEnd", "AWC", "TVN", "CTC"),
    @(7, "5:28 AM", 202.2, 232.2, 11.26, 1.29, 0),
    @(189, "7:00 AM", 5668.34, 5698.34, 0, 1, 0),
    @(309, "7:59 AM", 9248.33, 9278.33, 0, 1.065, 0),
    @(335, "8:12 AM", 10035.655, 10065.655, 4.215, 1.55, 0),
    @(380, "8:35 AM", 11373.16, 11403.16, 2.33, 3.355, 0),
    @(384, "8:37 AM", 11516.145, 11546.145, 3.955, 28.035, 0),
    @(388, "8:39 AM", 11611.81, 11641.81, 0, 59.97, 0),
    @(391, "8:41 AM", 11729.405, 11759.405, 0, 21.885, 0),
    @(395, "8:42 AM", 11822.4, 11852.4, 0, 8.73, 0),
    @(398, "8:44 AM", 11919.3225, 11949.3225, 1, 2.18, 0.25),
    @(402, "8:46 AM", 12035.575, 12065.575, 0, 2.95, 0),
    @(406, "8:48 AM", 12167.005, 12197.005, 0, 9.125, 0),
    @(425, "8:57 AM", 12725.265, 12755.265, 4.76, 43.895, 0),
    @(428, "8:59 AM", 12831.085, 12861.085, 0, 51.785, 0),
    @(432, "9:01 AM", 12941.465, 12971.465, 2.38, 48.73, 0.5),
    @(436, "9:03 AM", 13056.61, 13086.61, 3.23, 13.92, 0),
    @(439, "9:05 AM", 13169.08, 13199.08, 0, 1.32, 0),
    @(443, "9:06 AM", 13276.226667, 13306.226667, 1.58, 5.873333, 0.333333),
    @(446, "9:08 AM", 13376.97, 13406.97, 0, 18.32, 0),
    @(449, "9:10 AM", 13469.65, 13499.65, 0, 3.12, 0),
    @(453, "9:11 AM", 13570.62, 13600.62, 0, 21, 0),
    @(461, "9:15 AM", 13803.455, 13833.455, 0, 45.39, 0),
    @(465, "9:17 AM", 13937.845, 13967.845, 0, 30.79, 0),
    @(470, "9:20 AM", 14073.82, 14103.82, 0, 30.935, 0),
    @(473, "9:22 AM", 14186.425, 14216.425, 0, 4.27, 0),
    @(478, "9:24 AM", 14324.6875, 14354.6875, 0, 11.85, 0),
    @(482, "9:26 AM", 14442.785, 14472.785, 8.75, 7.99, 1),
    @(486, "9:28 AM", 14565.195, 14595.195, 0, 27.295, 0),
    @(496, "9:33 AM", 14854.345, 14884.345, 0, 19.64, 0),
    @(499, "9:34 AM", 14960.895, 14990.895, 0, 18.025, 0),
    @(504, "9:37 AM", 15106.36, 15136.36, 0, 52.755, 0),
    @(510, "9:40 AM", 15282.67, 15312.67, 0, 24.24, 0),
    @(518, "9:44 AM", 15517.59, 15547.59, 20.516667, 0.58, 0),
    @(568, "10:09 A", 17019.785, 17049.785, 8.845000000000001, 1.025, 1),
    @(576, "10:13 A", 17250.1, 17280.1, 0, 3.83, 0),
    @(581, "10:15 A", 17413.61, 17443.61, 0, 676.14, 0),
    @(983, "1:36 PM", 29464.66, 29494.66, 0.803333, 0.78, 0),
    @(999, "1:44 PM", 29954.343333, 29984.343333, 3.856667, 0.826667, 0.666667),
    @(1043, "2:06 PM", 31284.45, 31314.45, 0, 1.12, 0),
    @(1148, "2:59 PM", 34414.095, 34444.095, 0, 0.8, 0),
    @(1383, "4:56 PM", 41472.38, 41502.38, 7.295, 1.325, 0),
    @(1387, "4:58 PM", 41596.02, 41626.02, 0, 1.94, 0),
    @(1397, "5:04 PM", 41905.92, 41935.92, 7.78, 1, 1),
    @(1458, "5:34 PM", 43728.35, 43758.35, 3.405, 0.5649999999999999, 0),
    @(1547, "6:19 PM", 46409.81, 46439.81, 0, 1.36, 0),
    @(1551, "6:20 PM", 46502.28, 46532.28, 39.27, 1.01, 1),
    @(1558, "6:24 PM", 46717.07, 46747.07, 0, 0.64, 0),
    @(1561, "6:25 PM", 46813.12, 46843.12, 0, 1.03, 0),
    @(1608, "6:49 PM", 48231.23, 48261.23, 0, 0.695, 0)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Apply the same header formatting (bold, bordered, centered) used for A1:E1 to the new F1 and G1 headers
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 51 held the previous last data row (index 1608); that data now lives in row 50, so clear the old leftover row.
$ws.Rows.Item(51).ClearContents()
